# Auto-generated PowerShell COM-interop script
# Applies the 07:17:59 scrape update to the LÍNEA 141 schedule workbook
# (sheets: LP1912, LP1912-215, 6203-6173)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = 'Última actualización: 07:17:59'
$ws.Range("A3").Value = 'Total filas: 81'
$data = New-Object 'object[,]' 81,5
$data[0,0] = '04:45:48'
$data[0,1] = '04:45'
$data[0,2] = '215A_EL PATO'
$data[0,3] = 0
$data[0,4] = 'LP1912'
$data[1,0] = '04:21:09'
$data[1,1] = '04:46'
$data[1,2] = '215A_EL PATO'
$data[1,3] = 25
$data[1,4] = 'LP1912'
$data[2,0] = '04:45:48'
$data[2,1] = '04:53'
$data[2,2] = '11_ETCHEVERRY'
$data[2,3] = 8
$data[2,4] = 'LP1912'
$data[3,0] = '04:56:30'
$data[3,1] = '05:16'
$data[3,2] = '17_ROMERO'
$data[3,3] = 20
$data[3,4] = 'LP1912'
$data[4,0] = '04:56:30'
$data[4,1] = '05:22'
$data[4,2] = '23_HERNANDEZ'
$data[4,3] = 26
$data[4,4] = 'LP1912'
$data[5,0] = '05:24:16'
$data[5,1] = '05:25'
$data[5,2] = '23_HERNANDEZ'
$data[5,3] = 1
$data[5,4] = 'LP1912'
$data[6,0] = '04:56:30'
$data[6,1] = '05:34'
$data[6,2] = '215B_EL PATO'
$data[6,3] = 38
$data[6,4] = 'LP1912'
$data[7,0] = '05:24:16'
$data[7,1] = '05:35'
$data[7,2] = '215B_EL PATO'
$data[7,3] = 11
$data[7,4] = 'LP1912'
$data[8,0] = '05:24:16'
$data[8,1] = '05:46'
$data[8,2] = '15_ABASTO'
$data[8,3] = 22
$data[8,4] = 'LP1912'
$data[9,0] = '05:24:16'
$data[9,1] = '05:54'
$data[9,2] = '10_OLMOS'
$data[9,3] = 30
$data[9,4] = 'LP1912'
$data[10,0] = '05:55:02'
$data[10,1] = '05:55'
$data[10,2] = '10_OLMOS'
$data[10,3] = 0
$data[10,4] = 'LP1912'
$data[11,0] = '05:24:16'
$data[11,1] = '06:04'
$data[11,2] = '16_SANTA ANA'
$data[11,3] = 40
$data[11,4] = 'LP1912'
$data[12,0] = '05:24:16'
$data[12,1] = '06:11'
$data[12,2] = '215A_EL PATO'
$data[12,3] = 47
$data[12,4] = 'LP1912'
$data[13,0] = '05:55:02'
$data[13,1] = '06:12'
$data[13,2] = '215A_EL PATO'
$data[13,3] = 17
$data[13,4] = 'LP1912'
$data[14,0] = '05:55:02'
$data[14,1] = '06:14'
$data[14,2] = '225_HARAS DEL SUR'
$data[14,3] = 19
$data[14,4] = 'LP1912'
$data[15,0] = '04:56:30'
$data[15,1] = '06:18'
$data[15,2] = '16_SANTA ANA'
$data[15,3] = 82
$data[15,4] = 'LP1912'
$data[16,0] = '05:55:02'
$data[16,1] = '06:21'
$data[16,2] = '26_HERNANDEZ'
$data[16,3] = 26
$data[16,4] = 'LP1912'
$data[17,0] = '04:45:48'
$data[17,1] = '06:24'
$data[17,2] = '16_SANTA ANA'
$data[17,3] = 99
$data[17,4] = 'LP1912'
$data[18,0] = '05:55:02'
$data[18,1] = '06:27'
$data[18,2] = '23_HERNANDEZ'
$data[18,3] = 32
$data[18,4] = 'LP1912'
$data[19,0] = '06:25:28'
$data[19,1] = '06:29'
$data[19,2] = '86_EST CHICA-ESC AGRARIA'
$data[19,3] = 4
$data[19,4] = 'LP1912'
$data[20,0] = '06:25:28'
$data[20,1] = '06:30'
$data[20,2] = '23_HERNANDEZ'
$data[20,3] = 5
$data[20,4] = 'LP1912'
$data[21,0] = '05:55:02'
$data[21,1] = '06:30'
$data[21,2] = '86_EST CHICA-ESC AGRARIA'
$data[21,3] = 35
$data[21,4] = 'LP1912'
$data[22,0] = '05:55:02'
$data[22,1] = '06:31'
$data[22,2] = '16_SANTA ANA'
$data[22,3] = 36
$data[22,4] = 'LP1912'
$data[23,0] = '06:25:28'
$data[23,1] = '06:44'
$data[23,2] = '225_C ROCA-H SUR'
$data[23,3] = 19
$data[23,4] = 'LP1912'
$data[24,0] = '06:25:28'
$data[24,1] = '06:46'
$data[24,2] = '215C_EL PATO'
$data[24,3] = 21
$data[24,4] = 'LP1912'
$data[25,0] = '05:55:02'
$data[25,1] = '06:47'
$data[25,2] = '215C_EL PATO'
$data[25,3] = 52
$data[25,4] = 'LP1912'
$data[26,0] = '06:54:06'
$data[26,1] = '06:55'
$data[26,2] = '215C_EL PATO'
$data[26,3] = 1
$data[26,4] = 'LP1912'
$data[27,0] = '06:54:06'
$data[27,1] = '06:55'
$data[27,2] = '14_ABASTO'
$data[27,3] = 1
$data[27,4] = 'LP1912'
$data[28,0] = '06:25:28'
$data[28,1] = '06:59'
$data[28,2] = '14_ABASTO'
$data[28,3] = 34
$data[28,4] = 'LP1912'
$data[29,0] = '05:55:02'
$data[29,1] = '07:00'
$data[29,2] = '14_ABASTO'
$data[29,3] = 65
$data[29,4] = 'LP1912'
$data[30,0] = '06:54:06'
$data[30,1] = '07:01'
$data[30,2] = '16_SANTA ANA'
$data[30,3] = 7
$data[30,4] = 'LP1912'
$data[31,0] = '06:25:28'
$data[31,1] = '07:05'
$data[31,2] = '15_ABASTO'
$data[31,3] = 40
$data[31,4] = 'LP1912'
$data[32,0] = '06:54:06'
$data[32,1] = '07:05'
$data[32,2] = '23_HERNANDEZ'
$data[32,3] = 11
$data[32,4] = 'LP1912'
$data[33,0] = '06:54:06'
$data[33,1] = '07:07'
$data[33,2] = '15_ABASTO'
$data[33,3] = 13
$data[33,4] = 'LP1912'
$data[34,0] = '06:54:06'
$data[34,1] = '07:07'
$data[34,2] = '225_GOMEZ'
$data[34,3] = 13
$data[34,4] = 'LP1912'
$data[35,0] = '06:25:28'
$data[35,1] = '07:11'
$data[35,2] = '215A_EL PATO'
$data[35,3] = 46
$data[35,4] = 'LP1912'
$data[36,0] = '06:54:06'
$data[36,1] = '07:12'
$data[36,2] = '215A_EL PATO'
$data[36,3] = 18
$data[36,4] = 'LP1912'
$data[37,0] = '06:25:28'
$data[37,1] = '07:15'
$data[37,2] = '11_ETCHEVERRY'
$data[37,3] = 50
$data[37,4] = 'LP1912'
$data[38,0] = '06:54:06'
$data[38,1] = '07:16'
$data[38,2] = '11_ETCHEVERRY'
$data[38,3] = 22
$data[38,4] = 'LP1912'
$data[39,0] = '06:54:06'
$data[39,1] = '07:17'
$data[39,2] = '16_SANTA ANA'
$data[39,3] = 23
$data[39,4] = 'LP1912'
$data[40,0] = '07:17:59'
$data[40,1] = '07:20'
$data[40,2] = '26_HERNANDEZ'
$data[40,3] = 3
$data[40,4] = 'LP1912'
$data[41,0] = '06:54:06'
$data[41,1] = '07:21'
$data[41,2] = '26_HERNANDEZ'
$data[41,3] = 27
$data[41,4] = 'LP1912'
$data[42,0] = '06:54:06'
$data[42,1] = '07:23'
$data[42,2] = '10_OLMOS'
$data[42,3] = 29
$data[42,4] = 'LP1912'
$data[43,0] = '07:17:59'
$data[43,1] = '07:31'
$data[43,2] = '11_ETCHEVERRY'
$data[43,3] = 14
$data[43,4] = 'LP1912'
$data[44,0] = '07:17:59'
$data[44,1] = '07:31'
$data[44,2] = '16_SANTA ANA'
$data[44,3] = 14
$data[44,4] = 'LP1912'
$data[45,0] = '07:17:59'
$data[45,1] = '07:31'
$data[45,2] = '84_COLONIA URQUIZA-ESC 49'
$data[45,3] = 14
$data[45,4] = 'LP1912'
$data[46,0] = '06:54:06'
$data[46,1] = '07:32'
$data[46,2] = '11_ETCHEVERRY'
$data[46,3] = 38
$data[46,4] = 'LP1912'
$data[47,0] = '06:54:06'
$data[47,1] = '07:32'
$data[47,2] = '84_COLONIA URQUIZA-ESC 49'
$data[47,3] = 38
$data[47,4] = 'LP1912'
$data[48,0] = '05:55:02'
$data[48,1] = '07:32'
$data[48,2] = '16_SANTA ANA'
$data[48,3] = 97
$data[48,4] = 'LP1912'
$data[49,0] = '07:17:59'
$data[49,1] = '07:35'
$data[49,2] = '23_HERNANDEZ'
$data[49,3] = 18
$data[49,4] = 'LP1912'
$data[50,0] = '07:17:59'
$data[50,1] = '07:36'
$data[50,2] = '27_EL RETIRO'
$data[50,3] = 19
$data[50,4] = 'LP1912'
$data[51,0] = '06:54:06'
$data[51,1] = '07:37'
$data[51,2] = '27_EL RETIRO'
$data[51,3] = 43
$data[51,4] = 'LP1912'
$data[52,0] = '07:17:59'
$data[52,1] = '07:38'
$data[52,2] = '10_OLMOS'
$data[52,3] = 21
$data[52,4] = 'LP1912'
$data[53,0] = '06:54:06'
$data[53,1] = '07:39'
$data[53,2] = '10_OLMOS'
$data[53,3] = 45
$data[53,4] = 'LP1912'
$data[54,0] = '07:17:59'
$data[54,1] = '07:46'
$data[54,2] = '16_SANTA ANA'
$data[54,3] = 29
$data[54,4] = 'LP1912'
$data[55,0] = '07:17:59'
$data[55,1] = '07:47'
$data[55,2] = '14_ABASTO'
$data[55,3] = 30
$data[55,4] = 'LP1912'
$data[56,0] = '06:54:06'
$data[56,1] = '07:48'
$data[56,2] = '14_ABASTO'
$data[56,3] = 54
$data[56,4] = 'LP1912'
$data[57,0] = '07:17:59'
$data[57,1] = '07:51'
$data[57,2] = '215D_EL PATO'
$data[57,3] = 34
$data[57,4] = 'LP1912'
$data[58,0] = '06:54:06'
$data[58,1] = '07:52'
$data[58,2] = '215D_EL PATO'
$data[58,3] = 58
$data[58,4] = 'LP1912'
$data[59,0] = '07:17:59'
$data[59,1] = '07:59'
$data[59,2] = '23_HERNANDEZ'
$data[59,3] = 42
$data[59,4] = 'LP1912'
$data[60,0] = '06:25:28'
$data[60,1] = '08:03'
$data[60,2] = '23_HERNANDEZ'
$data[60,3] = 98
$data[60,4] = 'LP1912'
$data[61,0] = '07:17:59'
$data[61,1] = '08:03'
$data[61,2] = '11_ETCHEVERRY'
$data[61,3] = 46
$data[61,4] = 'LP1912'
$data[62,0] = '06:54:06'
$data[62,1] = '08:06'
$data[62,2] = '23_HERNANDEZ'
$data[62,3] = 72
$data[62,4] = 'LP1912'
$data[63,0] = '07:17:59'
$data[63,1] = '08:11'
$data[63,2] = '15_ABASTO'
$data[63,3] = 54
$data[63,4] = 'LP1912'
$data[64,0] = '07:17:59'
$data[64,1] = '08:12'
$data[64,2] = '10_OLMOS'
$data[64,3] = 55
$data[64,4] = 'LP1912'
$data[65,0] = '06:54:06'
$data[65,1] = '08:12'
$data[65,2] = '15_ABASTO'
$data[65,3] = 78
$data[65,4] = 'LP1912'
$data[66,0] = '07:17:59'
$data[66,1] = '08:20'
$data[66,2] = '26_HERNANDEZ'
$data[66,3] = 63
$data[66,4] = 'LP1912'
$data[67,0] = '06:54:06'
$data[67,1] = '08:21'
$data[67,2] = '26_HERNANDEZ'
$data[67,3] = 87
$data[67,4] = 'LP1912'
$data[68,0] = '07:17:59'
$data[68,1] = '08:22'
$data[68,2] = '16_P MOR-SANTA ANA'
$data[68,3] = 65
$data[68,4] = 'LP1912'
$data[69,0] = '07:17:59'
$data[69,1] = '08:22'
$data[69,2] = '215B_EL PATO'
$data[69,3] = 65
$data[69,4] = 'LP1912'
$data[70,0] = '06:54:06'
$data[70,1] = '08:23'
$data[70,2] = '215B_EL PATO'
$data[70,3] = 89
$data[70,4] = 'LP1912'
$data[71,0] = '06:54:06'
$data[71,1] = '08:23'
$data[71,2] = '16_P MOR-SANTA ANA'
$data[71,3] = 89
$data[71,4] = 'LP1912'
$data[72,0] = '07:17:59'
$data[72,1] = '08:26'
$data[72,2] = '84_COLONIA URQUIZA-ESC 49'
$data[72,3] = 69
$data[72,4] = 'LP1912'
$data[73,0] = '06:54:06'
$data[73,1] = '08:27'
$data[73,2] = '84_COLONIA URQUIZA-ESC 49'
$data[73,3] = 93
$data[73,4] = 'LP1912'
$data[74,0] = '07:17:59'
$data[74,1] = '08:41'
$data[74,2] = '81_EL PELIGRO'
$data[74,3] = 84
$data[74,4] = 'LP1912'
$data[75,0] = '06:54:06'
$data[75,1] = '08:42'
$data[75,2] = '81_EL PELIGRO'
$data[75,3] = 108
$data[75,4] = 'LP1912'
$data[76,0] = '07:17:59'
$data[76,1] = '08:43'
$data[76,2] = '14_ABASTO'
$data[76,3] = 86
$data[76,4] = 'LP1912'
$data[77,0] = '07:17:59'
$data[77,1] = '08:53'
$data[77,2] = '17_ROMERO'
$data[77,3] = 96
$data[77,4] = 'LP1912'
$data[78,0] = '07:17:59'
$data[78,1] = '09:01'
$data[78,2] = '215A_EL PATO'
$data[78,3] = 104
$data[78,4] = 'LP1912'
$data[79,0] = '07:17:59'
$data[79,1] = '09:10'
$data[79,2] = '16_P MOR-SANTA ANA'
$data[79,3] = 113
$data[79,4] = 'LP1912'
$data[80,0] = '07:17:59'
$data[80,1] = '09:16'
$data[80,2] = '27_EL RETIRO'
$data[80,3] = 119
$data[80,4] = 'LP1912'
$ws.Range("A6:E86").Value = $data

$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = 'Última actualización: 07:17:59'
$ws.Range("A3").Value = 'Total filas: 16'
$data = New-Object 'object[,]' 16,5
$data[0,0] = '04:45:48'
$data[0,1] = '04:45'
$data[0,2] = '215A_EL PATO'
$data[0,3] = 0
$data[0,4] = 'LP1912'
$data[1,0] = '04:21:09'
$data[1,1] = '04:46'
$data[1,2] = '215A_EL PATO'
$data[1,3] = 25
$data[1,4] = 'LP1912'
$data[2,0] = '04:56:30'
$data[2,1] = '05:34'
$data[2,2] = '215B_EL PATO'
$data[2,3] = 38
$data[2,4] = 'LP1912'
$data[3,0] = '05:24:16'
$data[3,1] = '05:35'
$data[3,2] = '215B_EL PATO'
$data[3,3] = 11
$data[3,4] = 'LP1912'
$data[4,0] = '05:24:16'
$data[4,1] = '06:11'
$data[4,2] = '215A_EL PATO'
$data[4,3] = 47
$data[4,4] = 'LP1912'
$data[5,0] = '05:55:02'
$data[5,1] = '06:12'
$data[5,2] = '215A_EL PATO'
$data[5,3] = 17
$data[5,4] = 'LP1912'
$data[6,0] = '06:25:28'
$data[6,1] = '06:46'
$data[6,2] = '215C_EL PATO'
$data[6,3] = 21
$data[6,4] = 'LP1912'
$data[7,0] = '05:55:02'
$data[7,1] = '06:47'
$data[7,2] = '215C_EL PATO'
$data[7,3] = 52
$data[7,4] = 'LP1912'
$data[8,0] = '06:54:06'
$data[8,1] = '06:55'
$data[8,2] = '215C_EL PATO'
$data[8,3] = 1
$data[8,4] = 'LP1912'
$data[9,0] = '06:25:28'
$data[9,1] = '07:11'
$data[9,2] = '215A_EL PATO'
$data[9,3] = 46
$data[9,4] = 'LP1912'
$data[10,0] = '06:54:06'
$data[10,1] = '07:12'
$data[10,2] = '215A_EL PATO'
$data[10,3] = 18
$data[10,4] = 'LP1912'
$data[11,0] = '07:17:59'
$data[11,1] = '07:51'
$data[11,2] = '215D_EL PATO'
$data[11,3] = 34
$data[11,4] = 'LP1912'
$data[12,0] = '06:54:06'
$data[12,1] = '07:52'
$data[12,2] = '215D_EL PATO'
$data[12,3] = 58
$data[12,4] = 'LP1912'
$data[13,0] = '07:17:59'
$data[13,1] = '08:22'
$data[13,2] = '215B_EL PATO'
$data[13,3] = 65
$data[13,4] = 'LP1912'
$data[14,0] = '06:54:06'
$data[14,1] = '08:23'
$data[14,2] = '215B_EL PATO'
$data[14,3] = 89
$data[14,4] = 'LP1912'
$data[15,0] = '07:17:59'
$data[15,1] = '09:01'
$data[15,2] = '215A_EL PATO'
$data[15,3] = 104
$data[15,4] = 'LP1912'
$ws.Range("A6:E21").Value = $data

$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = 'Última actualización: 07:17:59'
$ws.Range("A3").Value = 'Total filas: 15'
$data = New-Object 'object[,]' 15,5
$data[0,0] = '04:45:48'
$data[0,1] = '05:43'
$data[0,2] = '215A_LA PLATA'
$data[0,3] = 58
$data[0,4] = 'L6173'
$data[1,0] = '05:24:16'
$data[1,1] = '05:44'
$data[1,2] = '215A_LA PLATA'
$data[1,3] = 20
$data[1,4] = 'L6173'
$data[2,0] = '04:45:48'
$data[2,1] = '06:08'
$data[2,2] = '215A_LA PLATA'
$data[2,3] = 83
$data[2,4] = 'L6173'
$data[3,0] = '05:55:02'
$data[3,1] = '06:09'
$data[3,2] = '215A_LA PLATA'
$data[3,3] = 14
$data[3,4] = 'L6173'
$data[4,0] = '04:45:48'
$data[4,1] = '06:32'
$data[4,2] = '215C_LA PLATA'
$data[4,3] = 107
$data[4,4] = 'L6203'
$data[5,0] = '06:25:28'
$data[5,1] = '06:33'
$data[5,2] = '215C_LA PLATA'
$data[5,3] = 8
$data[5,4] = 'L6203'
$data[6,0] = '06:54:06'
$data[6,1] = '07:00'
$data[6,2] = '215B_LP-P MOR-1 Y 57'
$data[6,3] = 6
$data[6,4] = 'L6173'
$data[7,0] = '07:17:59'
$data[7,1] = '07:34'
$data[7,2] = '215A_LA PLATA'
$data[7,3] = 17
$data[7,4] = 'L6173'
$data[8,0] = '06:54:06'
$data[8,1] = '07:35'
$data[8,2] = '215A_LA PLATA'
$data[8,3] = 41
$data[8,4] = 'L6173'
$data[9,0] = '06:25:28'
$data[9,1] = '08:07'
$data[9,2] = '215C_LA PLATA'
$data[9,3] = 102
$data[9,4] = 'L6203'
$data[10,0] = '06:54:06'
$data[10,1] = '08:13'
$data[10,2] = '215C_LA PLATA'
$data[10,3] = 79
$data[10,4] = 'L6203'
$data[11,0] = '07:17:59'
$data[11,1] = '08:19'
$data[11,2] = '215C_LA PLATA'
$data[11,3] = 62
$data[11,4] = 'L6203'
$data[12,0] = '07:17:59'
$data[12,1] = '08:34'
$data[12,2] = '215A_LA PLATA'
$data[12,3] = 77
$data[12,4] = 'L6173'
$data[13,0] = '06:54:06'
$data[13,1] = '08:35'
$data[13,2] = '215A_LA PLATA'
$data[13,3] = 101
$data[13,4] = 'L6173'
$data[14,0] = '07:17:59'
$data[14,1] = '09:08'
$data[14,2] = '215D_LA PLATA'
$data[14,3] = 111
$data[14,4] = 'L6203'
$ws.Range("A6:E20").Value = $data
